# Apply crypto price/volume updates per commit diff (Mon Feb 26 09:25:05 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "51.158.35"
$ws.Cells.Item(2, 5).Value = "  -0.80%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.066.98"
$ws.Cells.Item(3, 5).Value = "  +1.68%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.18%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "386.63"
$ws.Cells.Item(5, 5).Value = "  +1.95%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "102.51"
$ws.Cells.Item(6, 5).Value = "  +0.12%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.536"
$ws.Cells.Item(7, 5).Value = "  -1.62%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.04%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.580"
$ws.Cells.Item(9, 5).Value = "  -1.49%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "36.68"
$ws.Cells.Item(10, 5).Value = "  -0.16%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +0.18%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0849"
$ws.Cells.Item(12, 5).Value = "  -1.39%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "3.576.31"
$ws.Cells.Item(13, 5).Value = "  +2.32%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "18.32"
$ws.Cells.Item(14, 5).Value = "  -0.43%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "7.68"
$ws.Cells.Item(15, 5).Value = "  -0.24%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "3.090.66"
$ws.Cells.Item(16, 5).Value = "  +2.18%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.986"
$ws.Cells.Item(17, 5).Value = "  +1.52%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "10.60"
$ws.Cells.Item(18, 5).Value = "  -0.07%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "51.276.82"
$ws.Cells.Item(19, 5).Value = "  -0.55%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "3.20"
$ws.Cells.Item(20, 5).Value = "  +4.31%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "12.37"
$ws.Cells.Item(21, 5).Value = "  -0.28%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.0₃0957"
$ws.Cells.Item(22, 5).Value = "  -0.48%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "69.78"
$ws.Cells.Item(23, 5).Value = "  -0.18%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "264.13"
$ws.Cells.Item(24, 5).Value = "  -0.80%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "3.14"
$ws.Cells.Item(25, 5).Value = "  -0.23%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "7.89"
$ws.Cells.Item(26, 5).Value = "  -4.37%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "26.96"
$ws.Cells.Item(27, 5).Value = "  +3.14%  "

# Row 28
$ws.Cells.Item(28, 2).Value = "Dai"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.00"
$ws.Cells.Item(28, 5).Value = "  +0.01%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "RenderToken"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "7.18"
$ws.Cells.Item(29, 5).Value = "  -5.48%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.164"
$ws.Cells.Item(30, 5).Value = "  -4.60%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.105"
$ws.Cells.Item(31, 5).Value = "  -2.38%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "10.46"
$ws.Cells.Item(32, 5).Value = "  +2.14%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "35.51"
$ws.Cells.Item(33, 5).Value = "  +5.22%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.0472"
$ws.Cells.Item(34, 5).Value = "  +4.85%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.06"
$ws.Cells.Item(35, 5).Value = "  +0.09%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "49.73"
$ws.Cells.Item(36, 5).Value = "  -1.58%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +0.14%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "3.36"
$ws.Cells.Item(38, 5).Value = "  +2.22%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.288"
$ws.Cells.Item(39, 5).Value = "  -0.61%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "128.98"
$ws.Cells.Item(40, 5).Value = "  +4.22%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "ARBITRUM"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.83"
$ws.Cells.Item(41, 5).Value = "  -1.15%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "Celestia"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "16.48"
$ws.Cells.Item(42, 5).Value = "  -2.63%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "Stellar"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.115"
$ws.Cells.Item(43, 5).Value = "  -0.53%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "NEARProtocol"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "3.81"
$ws.Cells.Item(44, 5).Value = "  +1.37%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.47"
$ws.Cells.Item(45, 5).Value = "  -2.67%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "21.71"
$ws.Cells.Item(46, 5).Value = "  +0.37%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.50"
$ws.Cells.Item(47, 5).Value = "  +4.98%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.09"
$ws.Cells.Item(48, 5).Value = "  -0.03%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.059.40"
$ws.Cells.Item(49, 5).Value = "  +1.67%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "FraxShare"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "9.64"
$ws.Cells.Item(50, 5).Value = "  +14.27%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "Mantle"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.930"
$ws.Cells.Item(51, 5).Value = "  +18.44%  "
